$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.444.86'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.567.99'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").Formula = '="288.16"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("D7").Formula = '="0.3719"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("D8").Formula = '="48.23"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  -4.10%  '
$ws.Range("D9").Formula = '="0.3313"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  -2.14%  '
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("D11").Formula = '="0.07470"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Formula = '="20.61"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  -2.75%  '
$ws.Range("D14").Formula = '="5.929"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  -1.61%  '
$ws.Range("D15").Formula = '="6.904"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  -1.33%  '
$ws.Range("D16").Value = '1.567.00'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Formula = '="87.70"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -2.92%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Formula = '="0.06743"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").Formula = '="6.350"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").Formula = '="16.46"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("D24").Value = '22.446.87'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Formula = '="2.390"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("D26").Formula = '="2.563"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -3.32%  '
$ws.Range("D27").Formula = '="152.84"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +2.19%  '
$ws.Range("D28").Formula = '="19.67"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("D29").Formula = '="5.024"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("D30").Formula = '="124.12"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -0.52%  '
$ws.Range("D31").Value = '1.742.23'
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("E32").Value = '  -0.97%  '
$ws.Range("D33").Formula = '="2.010"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("E34").Value = '  -1.87%  '
$ws.Range("D35").Formula = '="9.750"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("D36").Formula = '="0.08318"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -0.77%  '
$ws.Range("D37").Formula = '="0.02459"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  -0.61%  '
$ws.Range("E38").Value = '  -1.16%  '
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("D40").Formula = '="5.371"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("D41").Formula = '="1.286"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -4.60%  '
$ws.Range("D42").Formula = '="11.27"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("D44").Formula = '="1.001"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Formula = '="13.84"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -1.62%  '
$ws.Range("D46").Formula = '="0.6145"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +4.71%  '
$ws.Range("D47").Formula = '="3.772"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("E48").Value = '  -0.78%  '
$ws.Range("D49").Formula = '="125.71"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("E50").Value = '  -2.13%  '
$ws.Range("D51").Formula = '="0.07218"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -1.29%  '
$excel.CutCopyMode = $false
